$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume update (GitHub Actions data refresh)
# D-column (Price) values are numeric-looking text (e.g. "569.50", "0.0000233")
# that Excel would otherwise auto-coerce to a Number and normalize (dropping
# trailing zeros / losing the original text type). Force the cell to Text
# format before assigning, then restore the default "Normal" style so no
# extra style index is left attached to the cell (matches source formatting).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.484.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.889.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.41%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.888.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.76%  "
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.02%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.367.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.494.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.890.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("E31").Value = "  -5.48%  "
$ws.Range("E32").Value = "  -7.77%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.62%  "
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  -6.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.66%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.266"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.688.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "337.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.11%  "
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.59%  "
